$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2120.3623
$ws.Range("I17").Value = 294.14285
$ws.Range("J17").Value = 2585.2183
$ws.Range("K17").Value = 882.4285500000001
$ws.Range("L17").Value = 7755.6549
$ws.Range("M17").Value = -714.4285500000001
$ws.Range("N17").Value = -8091.6549
$ws.Range("H51").Value = 7716.3335
$ws.Range("J51").Value = 8066.6665
$ws.Range("L51").Value = 8066.6665
$ws.Range("N51").Value = -9034.666499999999
$ws.Range("H88").Value = 13080.417
$ws.Range("J88").Value = 17997.125
$ws.Range("L88").Value = 17997.125
$ws.Range("N88").Value = -18809.125
$ws.Range("H91").Value = 13080.417
$ws.Range("J91").Value = 17997.125
$ws.Range("L91").Value = 17997.125
$ws.Range("N91").Value = -20805.125
$ws.Range("H97").Value = 1421.75
$ws.Range("I97").Value = 1699
$ws.Range("K97").Value = 5097
$ws.Range("M97").Value = -4601
$ws.Range("H116").Value = 8478.799999999999
$ws.Range("I116").Value = 7112.4287
$ws.Range("K116").Value = 7112.4287
$ws.Range("M116").Value = -3670.4287
$ws.Range("H138").Value = 3118.9832
$ws.Range("I138").Value = 2075.318
$ws.Range("J138").Value = 3739.5405
$ws.Range("K138").Value = 6225.954000000001
$ws.Range("L138").Value = 11218.6215
$ws.Range("M138").Value = -1085.954000000001
$ws.Range("N138").Value = -21498.6215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 35.5
$ws.Range("I4").Value = 35.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 35.5
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 50016664
$ws.Range("J6").Value = 19998
$ws.Range("L6").Value = 19998
$ws.Range("N6").Value = -20344
$ws.Range("H32").Value = 4592.3486
$ws.Range("I32").Value = 4010.8809
$ws.Range("K32").Value = 4010.8809
$ws.Range("M32").Value = -3723.8809
$ws.Range("H61").Value = 6420.8125
$ws.Range("I61").Value = 3807
$ws.Range("K61").Value = 3807
$ws.Range("M61").Value = -3595
$ws.Range("H97").Value = 247.66667
$ws.Range("I97").Value = 288.08334
$ws.Range("K97").Value = 288.08334
$ws.Range("M97").Value = 207.91666
$ws.Range("H102").Value = 2048.3
$ws.Range("I102").Value = 2048.3
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2048.3
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 2992.08
$ws.Range("I132").Value = 2388.7
$ws.Range("K132").Value = 7166.099999999999
$ws.Range("M132").Value = -4636.099999999999
$ws.Range("H136").Value = 6420.8125
$ws.Range("I136").Value = 3807
$ws.Range("K136").Value = 11421
$ws.Range("M136").Value = -8871

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3376
$ws.Range("J11").Value = 4250
$ws.Range("L11").Value = 4250
$ws.Range("N11").Value = -4530
$ws.Range("H86").Value = 4774.1113
$ws.Range("J86").Value = 19007
$ws.Range("L86").Value = 19007
$ws.Range("N86").Value = -21253
$ws.Range("H89").Value = 4774.1113
$ws.Range("J89").Value = 19007
$ws.Range("L89").Value = 95035
$ws.Range("N89").Value = -106267

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 179.22223
$ws.Range("I7").Value = 73
$ws.Range("J7").Value = 285.44446
$ws.Range("K7").Value = 73
$ws.Range("L7").Value = 285.44446
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = -511.44446
$ws.Range("H9").Value = 133332
$ws.Range("J9").Value = 133332
$ws.Range("L9").Value = 133332
$ws.Range("N9").Value = -133668
$ws.Range("H105").Value = 11502.5
$ws.Range("I105").Value = 8999.666999999999
$ws.Range("J105").Value = 19011
$ws.Range("K105").Value = 8999.666999999999
$ws.Range("L105").Value = 19011
$ws.Range("M105").Value = -7252.666999999999
$ws.Range("N105").Value = -22505
$ws.Range("H141").Value = 132798.25
$ws.Range("J141").Value = 132798.25
$ws.Range("L141").Value = 132798.25
$ws.Range("N141").Value = -143158.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1902.091
$ws.Range("J86").Value = 2745.8572
$ws.Range("L86").Value = 8237.571599999999
$ws.Range("N86").Value = -10609.5716
$ws.Range("H89").Value = 1902.091
$ws.Range("J89").Value = 2745.8572
$ws.Range("L89").Value = 24712.7148
$ws.Range("N89").Value = -36568.7148
$ws.Range("H116").Value = 2674522.2
$ws.Range("I116").Value = 4004267.5
$ws.Range("J116").Value = 15032
$ws.Range("K116").Value = 12012802.5
$ws.Range("L116").Value = 45096
$ws.Range("M116").Value = -12009360.5
$ws.Range("N116").Value = -51980
$ws.Range("H136").Value = 2027
$ws.Range("I136").Value = 2031.5
$ws.Range("K136").Value = 6094.5
$ws.Range("M136").Value = -994.5
$ws.Range("H137").Value = 3039.75
$ws.Range("I137").Value = 1631.8889
$ws.Range("J137").Value = 4191.636
$ws.Range("K137").Value = 4895.6667
$ws.Range("L137").Value = 12574.908
$ws.Range("M137").Value = 204.3333000000002
$ws.Range("N137").Value = -22774.908
$ws.Range("H138").Value = 46020
$ws.Range("I138").Value = 56530
$ws.Range("J138").Value = 25000
$ws.Range("K138").Value = 169590
$ws.Range("L138").Value = 75000
$ws.Range("M138").Value = -164450
$ws.Range("N138").Value = -85280
$ws.Range("H139").Value = 3078.0833
$ws.Range("J139").Value = 19033
$ws.Range("L139").Value = 57099
$ws.Range("N139").Value = -67379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5177.4546
$ws.Range("J80").Value = 5768.75
$ws.Range("L80").Value = 5768.75
$ws.Range("N80").Value = -7764.75
$ws.Range("H83").Value = 5177.4546
$ws.Range("J83").Value = 5768.75
$ws.Range("L83").Value = 28843.75
$ws.Range("N83").Value = -38827.75
$ws.Range("H102").Value = 3374.75
$ws.Range("I102").Value = 3500
$ws.Range("K102").Value = 3500
$ws.Range("M102").Value = -1878
$ws.Range("H109").Value = 74789
$ws.Range("J109").Value = 74789
$ws.Range("L109").Value = 74789
$ws.Range("N109").Value = -76869
$ws.Range("H138").Value = 75428.336
$ws.Range("J138").Value = 75428.336
$ws.Range("L138").Value = 75428.336
$ws.Range("N138").Value = -85708.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4889
$ws.Range("I46").Value = 2666.6667
$ws.Range("J46").Value = 6000.1665
$ws.Range("K46").Value = 2666.6667
$ws.Range("L46").Value = 6000.1665
$ws.Range("M46").Value = -2478.6667
$ws.Range("N46").Value = -6376.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 49989
$ws.Range("J93").Value = 49989
$ws.Range("L93").Value = 49989
$ws.Range("N93").Value = -54981
$ws.Range("H107").Value = 557.6667
$ws.Range("I107").Value = 504
$ws.Range("J107").Value = 665
$ws.Range("K107").Value = 1512
$ws.Range("L107").Value = 1995
$ws.Range("M107").Value = 408
$ws.Range("N107").Value = -5835
$ws.Range("H113").Value = 686.73334
$ws.Range("I113").Value = 322.75
$ws.Range("K113").Value = 968.25
$ws.Range("M113").Value = 1201.75
$ws.Range("H136").Value = 1925
$ws.Range("I136").Value = 1436.2122
$ws.Range("K136").Value = 4308.6366
$ws.Range("M136").Value = -1758.6366
